# Second commit for practice
# Update a few test-data values on the "TestData" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

$ws.Range("C2").Value = "Serv014ice"
$ws.Range("H2").Value = "03062019"
$ws.Range("I2").Value = "10062019"

$ws.Range("I4").Select()
